$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert new column D ("Type"), shifting old D (Duration) -> E, old E (Comment) -> F
$ws.Columns.Item(4).Insert()

# 2. Set header for new Type column
$ws.Range("D1").Value2 = "Type"

# 3. Fill Type column (D) for existing rows
$ws.Range("D2").Value2  = "Research"
$ws.Range("D3").Value2  = "Tutorial"
$ws.Range("D4").Value2  = "Tutorial"
$ws.Range("D5").Value2  = "Tutorial"
$ws.Range("D6").Value2  = "Tutorial"
$ws.Range("D7").Value2  = "Tutorial"
$ws.Range("D8").Value2  = "Tutorial"
$ws.Range("D9").Value2  = "Tutorial"
$ws.Range("D10").Value2 = "Tutorial"
$ws.Range("D11").Value2 = "Experiment"
$ws.Range("D12").Value2 = "Documentation"
$ws.Range("D13").Value2 = "Documentation"
$ws.Range("D14").Value2 = "Documentation"
$ws.Range("D15").Value2 = "Documentation"
$ws.Range("D16").Value2 = "Research"

# 4. Update C16 task name
$ws.Range("C16").Value2 = "Existing Technologies"

# 5. Row 17: new task "Project Aim and Objectives"
$ws.Range("B17").Value2 = 0.91666666666666663
$ws.Range("C17").Value2 = "Project Aim and Objectives"
$ws.Range("D17").Value2 = "Documentation"
$ws.Range("E17").Value2 = 80
$ws.Range("F17").Value2 = "Research Gap, Aim and Objectives Paragraph, Literature Review: Sketch Strum vs Strumbar Task Plan"

# 6. Row 18: new task "Arduino Experiment 4X4 Matrix "
$ws.Range("C18").Value2 = "Arduino Experiment 4X4 Matrix "
$ws.Range("D18").Value2 = "Experiment"
$ws.Range("E18").Value2 = "TBA"
$ws.Range("F18").Value2 = "Create an Arduino 4X4 Numpad"

# Clear the old total formula that used to live at (old D19, now E19)
$ws.Range("E19").Clear()

# Row 19 stays blank like rows 17/18 in columns A and B (centered style retained)
$ws.Range("A19").HorizontalAlignment = -4108
$ws.Range("B19").HorizontalAlignment = -4108
$ws.Range("B19").NumberFormat = $ws.Range("B18").NumberFormat

# 7. Row 20: Minutes total
$ws.Range("D20").Value2 = "Minutes"
$ws.Range("E20").Formula = "=SUM(E2:E19)"

# 8. Row 21: Hours total
$ws.Range("D21").Value2 = "Hours"
$ws.Range("E21").Formula = "=E20 / 60"

# 9. Styling: header row gets bold + light-gray fill + centered text
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.Interior.ThemeColor = 2
$header.Interior.TintAndShade = -0.049989318521683403

# 10. Styling: date column (A2:A16) gets centered text + light-gray fill
$dateFill = $ws.Range("A2:A16")
$dateFill.HorizontalAlignment = -4108
$dateFill.Interior.ThemeColor = 2
$dateFill.Interior.TintAndShade = -0.049989318521683403

# 11. Styling: totals rows (D20:E21) bold
$totals = $ws.Range("D20:E21")
$totals.Font.Bold = $true

Write-Host "done content"
